$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new daily price record was added for "Feria Lagunitas de Puerto Montt -
# Zapallo italiano" on 2021-10-18. It belongs at the top of the data block
# (row 29, right after the header + first 27 data rows), so insert a fresh
# row there; Excel shifts every row below (old 29..133) down to 30..134,
# preserving their values untouched.
$ws.Rows.Item(29).Insert()

# Populate the newly-inserted row with the new record's data.
$ws.Range("A29").Value = 4
$ws.Range("B29").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C29").Value = "Los Lagos"
$ws.Range("D29").Value = 44487
$ws.Range("E29").Value = 10
$ws.Range("F29").Value = 100112032
$ws.Range("G29").Value = "Zapallo italiano"
$ws.Range("H29").Value = "Sin especificar"
$ws.Range("I29").Value = "Primera"
$ws.Range("J29").Value = 120
$ws.Range("K29").Value = 16000
$ws.Range("L29").Value = 16000
$ws.Range("M29").Value = 16000
$ws.Range("N29").Value = "`$/caja 50 unidades"
$ws.Range("O29").Value = "Región de Arica y Parinacota"
$ws.Range("P29").Value = 320
$ws.Range("Q29").Value = 50
$ws.Range("R29").Value = "Hortaliza"

# Match the date-number-format style used by the rest of column D.
$ws.Range("D29").NumberFormat = $ws.Range("D30").NumberFormat()
